$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 460, shifting existing rows 460..541 down to 461..542.
$ws.Rows.Item(460).Insert()

# Populate the newly inserted row 460 with the new weekly data point.
$ws.Cells.Item(460, 1).Value = 4
$ws.Cells.Item(460, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(460, 3).Value = "Los Lagos"
$ws.Cells.Item(460, 4).Value = 45209
$ws.Cells.Item(460, 5).Value = 10
$ws.Cells.Item(460, 6).Value = 100114014
$ws.Cells.Item(460, 7).Value = "Betarraga"
$ws.Cells.Item(460, 8).Value = "Sin especificar"
$ws.Cells.Item(460, 9).Value = "Primera"
$ws.Cells.Item(460, 10).Value = 1250
$ws.Cells.Item(460, 11).Value = 1000
$ws.Cells.Item(460, 12).Value = 1000
$ws.Cells.Item(460, 13).Value = 1000
$ws.Cells.Item(460, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(460, 15).Value = "Región Metropolitana"
$ws.Cells.Item(460, 16).Value = 200
$ws.Cells.Item(460, 17).Value = 5
$ws.Cells.Item(460, 18).Value = "Hortaliza"
